$wb = $excel.ActiveWorkbook

$sprint1 = $wb.Worksheets.Item("Sprint1")
$sprint4 = $wb.Worksheets.Item("Sprint4")

# --- Sprint4 (Sprint 4 sheet): add the E/F estimate/actual-hours columns ---
# Values mirror the Sprint1 sheet's E2:F9 pattern (75 / 60 per row).
$hours = @(75, 60)
foreach ($r in 2..9) {
    $e = $sprint4.Cells.Item($r, 5)
    $f = $sprint4.Cells.Item($r, 6)
    $e.Value = $hours[0]
    $f.Value = $hours[1]
    # Centered (both axes) formatting, matching the rest of the sheet's numeric columns.
    $e.HorizontalAlignment = -4108   # xlCenter
    $e.VerticalAlignment = -4108     # xlCenter
    $f.HorizontalAlignment = -4108
    $f.VerticalAlignment = -4108
}

# --- Selections on Sprint1 (mirrors the source range that was copied) ---
$sprint1.Activate()
$sprint1.Range("E2:F9").Select()

# --- Activate Sprint4 and select the pasted range; this also flips the
#     workbook's activeTab / each sheet's tabSelected+topLeftCell bookkeeping ---
$sprint4.Activate()
$sprint4.Range("E2:F9").Select()
